# Insert a new weekly price record as row 191 in the "Poroto verde" data
# sheet (Comercializadora del Agro de Limarí). Existing rows 191:220 shift
# down to 192:221, and the sheet's used range grows to A1:R221.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 191 and below down by one to make room for the new record.
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A191").Value = 2
$ws.Range("B191").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C191").Value = "Coquimbo"
$ws.Range("D191").Value = 44889
$ws.Range("E191").Value = 4
$ws.Range("F191").Value = 100112031
$ws.Range("G191").Value = "Poroto verde"
$ws.Range("H191").Value = "Magnum"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 400
$ws.Range("K191").Value = 25000
$ws.Range("L191").Value = 27000
$ws.Range("M191").Value = 26000
$ws.Range("N191").Value = "`$/caja 25 kilos"
$ws.Range("O191").Value = "Provincia de Limarí"
$ws.Range("P191").Value = 1040
$ws.Range("Q191").Value = 25
$ws.Range("R191").Value = "Hortaliza"
